$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '31.500.76'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +5.61%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.708.08'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +4.17%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.04%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '222.00'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +3.10%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.535'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +3.06%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.00'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.01%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '29.95'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +3.93%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.269'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +3.35%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0647'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +6.29%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0910'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +1.08%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.955.08'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +4.19%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.711.81'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +4.30%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.612'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +3.74%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '10.20'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +8.11%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '4.19'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +8.22%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '31.508.07'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +5.56%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '67.21'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +4.35%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '250.80'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +4.55%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.0₃0723'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +2.96%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.998'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.20%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '10.13'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +2.03%  '
$ws.Range("E23").Value = '  +2.82%  '
$ws.Range("E24").Value = '  -1.00%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '159.44'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +1.58%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '16.04'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +3.24%  '
$ws.Range("E27").Value = '  +3.17%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '6.78'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +2.40%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.00'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.04%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '3.86'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +14.11%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.0503'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +1.77%  '
$ws.Range("E32").Value = '  +4.05%  '
$ws.Range("E33").Value = '  +5.71%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.522.31'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +7.07%  '
$ws.Range("E35").Value = '  +2.80%  '
$ws.Range("E36").Value = '  +2.45%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '82.64'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +8.20%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.610'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +8.64%  '
$ws.Range("E39").Value = '  +4.38%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.73'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +0.35%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.32'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +0.85%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.853'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +2.54%  '
$ws.Range("E43").Value = '  +4.58%  '
$ws.Range("E44").Value = '  +0.84%  '
$ws.Range("E45").Value = '  +3.25%  '
$ws.Range("E46").Value = '  +0.03%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '52.12'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +5.62%  '
$ws.Range("E48").Value = '  +4.37%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.848.72'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +3.66%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0₆0119'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +9.90%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '93.56'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +0.42%  '
